# Update "想去人数" (column F) values across the four worksheets to the
# freshly scraped counts, as published at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Map of worksheet name -> list of (row, newValue) pairs for column F.
$updates = @{
    "展览" = @(
        @{Row=3;  Value=241},
        @{Row=4;  Value=828},
        @{Row=5;  Value=241},
        @{Row=6;  Value=410},
        @{Row=7;  Value=577},
        @{Row=10; Value=337},
        @{Row=11; Value=142},
        @{Row=12; Value=631},
        @{Row=13; Value=84},
        @{Row=14; Value=1789},
        @{Row=15; Value=336},
        @{Row=16; Value=2827},
        @{Row=17; Value=311},
        @{Row=18; Value=492},
        @{Row=19; Value=48},
        @{Row=20; Value=135}
    )
    "演出" = @(
        @{Row=4;  Value=45},
        @{Row=5;  Value=18},
        @{Row=6;  Value=118},
        @{Row=7;  Value=474},
        @{Row=15; Value=1}
    )
    "本地生活" = @(
        @{Row=2; Value=5308},
        @{Row=3; Value=313},
        @{Row=4; Value=235}
    )
    "全部类型" = @(
        @{Row=3;  Value=5308},
        @{Row=4;  Value=313},
        @{Row=6;  Value=235},
        @{Row=7;  Value=241},
        @{Row=9;  Value=45},
        @{Row=10; Value=18},
        @{Row=11; Value=118},
        @{Row=12; Value=474},
        @{Row=13; Value=828},
        @{Row=16; Value=241},
        @{Row=17; Value=410},
        @{Row=18; Value=577},
        @{Row=22; Value=337},
        @{Row=23; Value=142},
        @{Row=26; Value=631},
        @{Row=27; Value=84},
        @{Row=29; Value=1789},
        @{Row=30; Value=336},
        @{Row=31; Value=2829},
        @{Row=33; Value=311},
        @{Row=34; Value=492},
        @{Row=35; Value=48},
        @{Row=36; Value=1},
        @{Row=37; Value=135}
    )
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($entry in $updates[$sheetName]) {
        $ws.Cells.Item($entry.Row, 6).Value = $entry.Value
    }
}
